$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "intervention_type",
    "DRUG",
    "DRUG",
    "DEVICE",
    "DRUG",
    "PROCEDURE",
    "DRUG",
    "DRUG",
    "PROCEDURE",
    "DEVICE",
    "DRUG",
    "DRUG",
    "PROCEDURE",
    "OTHER",
    "PROCEDURE",
    "PROCEDURE",
    "PROCEDURE",
    "DEVICE",
    "DIETARY_SUPPLEMENT",
    "DRUG",
    "DRUG",
    "BIOLOGICAL",
    "DRUG",
    "DEVICE",
    "PROCEDURE",
    "DRUG",
    "DRUG",
    "GENETIC",
    "BEHAVIORAL",
    "PROCEDURE",
    "",
    "PROCEDURE",
    "GENETIC",
    "PROCEDURE",
    "DIETARY_SUPPLEMENT",
    "DEVICE",
    "OTHER",
    "DEVICE",
    "OTHER",
    "DIETARY_SUPPLEMENT",
    "DEVICE",
    "OTHER",
    "DRUG",
    "DRUG",
    "BIOLOGICAL",
    "DIETARY_SUPPLEMENT",
    "DRUG",
    "PROCEDURE",
    "DRUG",
    "BEHAVIORAL",
    "DEVICE",
    "DRUG",
    "PROCEDURE",
    "BIOLOGICAL",
    "OTHER",
    "DRUG",
    "DEVICE",
    "PROCEDURE",
    "DEVICE",
    "PROCEDURE",
    "OTHER",
    "",
    "GENETIC",
    "BIOLOGICAL",
    "PROCEDURE",
    "DRUG",
    "OTHER",
    "BIOLOGICAL",
    "BIOLOGICAL",
    "PROCEDURE",
    "DRUG",
    "GENETIC",
    "DEVICE",
    "DRUG",
    "PROCEDURE",
    "OTHER",
    "BEHAVIORAL",
    "DRUG",
    "PROCEDURE",
    "OTHER",
    "DEVICE",
    "DRUG",
    "DIAGNOSTIC_TEST",
    "PROCEDURE",
    "BIOLOGICAL",
    "DIETARY_SUPPLEMENT",
    "DEVICE",
    "PROCEDURE",
    "DEVICE",
    "BIOLOGICAL",
    "",
    "DRUG",
    "PROCEDURE",
    "",
    "DIAGNOSTIC_TEST",
    "RADIATION",
    "PROCEDURE",
    "PROCEDURE",
    "BEHAVIORAL",
    "OTHER",
    "OTHER",
    "PROCEDURE",
    "",
    "DIETARY_SUPPLEMENT",
    "PROCEDURE",
    "",
    "",
    "",
    ""
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $val = $values[$i]
    if ($val -ne "") {
        $ws.Cells.Item($row, 11).Value = $val
    }
}

# Apply the header style (bold, centered, bordered) from J1 to the new K1 header cell
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
